$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp text (row 1 header)
$ws.Range("A1").Value = "Datos actualizados a 17 de Octubre de 2020 a las 20:03"

# Row 4
$ws.Range("B4").Value = 8314560
$ws.Range("C4").Value = 26282
$ws.Range("D4").Value = 5406288
$ws.Range("E4").Value = 2684313
$ws.Range("G4").Value = 315
$ws.Range("H4").Value = 223959

# Row 5
$ws.Range("B5").Value = 7486714
$ws.Range("C5").Value = 56079
$ws.Range("D5").Value = 6587287
$ws.Range("E5").Value = 785648
$ws.Range("G5").Value = 747
$ws.Range("H5").Value = 113779

# Row 21
$ws.Range("B21").Value = 359655
$ws.Range("C21").Value = 2863
$ws.Range("D21").Value = 290000
$ws.Range("E21").Value = 59810
$ws.Range("G21").Value = 9
$ws.Range("H21").Value = 9845

# Row 33
$ws.Range("A33").Value = "Marruecos"
$ws.Range("B33").Value = 170911
$ws.Range("C33").Value = 3763
$ws.Range("D33").Value = 141381
$ws.Range("E33").Value = 26652
$ws.Range("G33").Value = 60
$ws.Range("H33").Value = 2878

# Row 34
$ws.Range("A34").Value = "Polonia"
$ws.Range("B34").Value = 167230
$ws.Range("C34").Value = 9622
$ws.Range("D34").Value = 90162
$ws.Range("E34").Value = 73544
$ws.Range("G34").Value = 84
$ws.Range("H34").Value = 3524

# Row 63
$ws.Range("A63").Value = "Libano"
$ws.Range("B63").Value = 61284
$ws.Range("C63").Value = 1171
$ws.Range("D63").Value = 27197
$ws.Range("E63").Value = 33570
$ws.Range("G63").Value = 8
$ws.Range("H63").Value = 517

# Row 64
$ws.Range("A64").Value = "Nigeria"
$ws.Range("B64").Value = 61194
$ws.Range("C64").Value = 0
$ws.Range("D64").Value = 52304
$ws.Range("E64").Value = 7771
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 1119

# Row 124
$ws.Range("B124").Value = 5765
$ws.Range("C124").Value = 19
$ws.Range("D124").Value = 5392
$ws.Range("E124").Value = 258
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 115

# Row 197
$ws.Range("B197").Value = 113
$ws.Range("C197").Value = 1
$ws.Range("D197").Value = 100
$ws.Range("E197").Value = 10
$ws.Range("G197").Value = 0
$ws.Range("H197").Value = 3

# Row 203
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("B203").Value = 33
$ws.Range("C203").Value = 1
$ws.Range("D203").Value = 27
$ws.Range("E203").Value = 6
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 204
$ws.Range("A204").Value = "Dominica"
$ws.Range("B204").Value = 33
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 29
$ws.Range("E204").Value = 4
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

# Row 205
$ws.Range("A205").Value = "Guam"
$ws.Range("B205").Value = 32
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 0
$ws.Range("E205").Value = 31
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 1

# Row 216
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 13
$ws.Range("E216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0

# Row 217
$ws.Range("A217").Value = "Montserrat"
$ws.Range("B217").Value = 13
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 12
$ws.Range("E217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 1
